$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment the "Folds" column (C2:C6) values by 1
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 4
$ws.Range("C6").Value = 5
